# QA Defect Report update:
#  - Retest/close defects 1-6 (Status -> Closed, add Reviewed-by + Comments)
#  - File a new defect #7 (Open) with description
#  - Add a "Reviewed by" / "Comments" column pair, colour-code rows by status
#  - Turn on AutoFilter over the table, filtered down to Status = Open,
#    freeze the header row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1. New columns F (Reviewed by) / G (Comments) on the header row
# ---------------------------------------------------------------------
$ws.Range("F2").Value = "Reviewed by"
$ws.Range("G2").Value = "Comments"

# ---------------------------------------------------------------------
# 2. Retest existing defects 1-6 : Status -> Closed, Reviewed by, Comments
# ---------------------------------------------------------------------
$ws.Range("D3").Value = "Closed"
$ws.Range("F3").Value = "Maksym"
$ws.Range("G3").Value = "Re-tested and passed. Closing defect. -VP"

$ws.Range("D4").Value = "Closed"
$ws.Range("F4").Value = "Maksym"
$ws.Range("G4").Value = "Re-tested and passed. Closing defect. -VP"

$ws.Range("D5").Value = "Closed"
$ws.Range("F5").Value = "Maksym"
$ws.Range("G5").Value = "changed url page generated to match naming convention -MM`nPost is displayed, but it is displayed as a response to an existing thread. New defect ID 7 created for this.  -VP"

$ws.Range("D6").Value = "Closed"
$ws.Range("F6").Value = "Maksym"
$ws.Range("G6").Value = "Never had this issue, pm if persists -MM`nCreated new defect ID 7 with more detailed description and steps to reproduce. VP"

$ws.Range("D7").Value = "Closed"
$ws.Range("F7").Value = "Maksym"
$ws.Range("G7").Value = "Re-tested and passed. Closing defect. -VP"

$ws.Range("D8").Value = "Closed"
$ws.Range("F8").Value = "Maksym"
$ws.Range("G8").Value = "Correction, only ' cause this error -MM`nRe-tested and passed. Closing defect. -VP"

# ---------------------------------------------------------------------
# 3. New defect #7 (row 9) - still Open
# ---------------------------------------------------------------------
$ws.Range("B9").Value = 7
$ws.Range("C9").Value = "New thread posts are displayed as a post response.`nSteps to reproduce:`nCreate a new Thread. `nAdd a Reply to the Thread.`nCreate another new Thread.`nView newly created thread.`nThe new thread appears as a Reply to the first Thread, e.g., the contents of the first thread are in the new thread."
$ws.Range("D9").Value = "Open"
$ws.Range("E9").Value = "vpersaud"

# ---------------------------------------------------------------------
# 4. Formatting : re-centre title, colour header + rows, wrap long text
# ---------------------------------------------------------------------
$ws.Range("A1:E1").HorizontalAlignment = -4108  # xlCenter, no fill

$headerRange = $ws.Range("A2:H2")
$headerRange.Font.Bold = $true
$headerRange.Interior.ThemeColor = 2
$headerRange.Interior.TintAndShade = -0.099978637043366805

$closedRows = @(3, 4, 5, 6, 7, 8)
foreach ($r in $closedRows) {
    $rowRange = $ws.Range("A$r`:G$r")
    $rowRange.Interior.ThemeColor = 9
    $rowRange.Interior.TintAndShade = 0.59999389629810485
    $ws.Range("C$r").WrapText = $true
    $ws.Range("G$r").WrapText = $true
}

$openRowRange = $ws.Range("A9:G9")
$openRowRange.Interior.ThemeColor = 7
$openRowRange.Interior.TintAndShade = 0.59999389629810485
$ws.Range("C9").WrapText = $true
$ws.Range("G9").WrapText = $true

# ---------------------------------------------------------------------
# 5. Row heights (Excel auto-grows wrapped rows; set explicitly here)
# ---------------------------------------------------------------------
$ws.Rows.Item(3).RowHeight = 29.5
$ws.Rows.Item(4).RowHeight = 29.5
$ws.Rows.Item(5).RowHeight = 73.75
$ws.Rows.Item(6).RowHeight = 59
$ws.Rows.Item(7).RowHeight = 29.5
$ws.Rows.Item(8).RowHeight = 176.5
$ws.Rows.Item(9).RowHeight = 132.75

# ---------------------------------------------------------------------
# 6. Column width tweaks
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 20.54296875
$ws.Columns.Item(4).ColumnWidth = 5.86328125
$ws.Columns.Item(6).ColumnWidth = 21.1328125
$ws.Columns.Item(7).ColumnWidth = 47.7265625

# ---------------------------------------------------------------------
# 7. AutoFilter the table, then filter Status (col D) down to "Open"
# ---------------------------------------------------------------------
$ws.Range("A2:H9").AutoFilter(4, "Open")

# ---------------------------------------------------------------------
# 8. Freeze the header row, select a cell below it
# ---------------------------------------------------------------------
$ws.Range("A2").Select()
$excel.ActiveWindow.FreezePanes = $true
$ws.Range("F16").Select()

Write-Output "QA defect report updated"
